$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245, shifting existing rows 245:321 down to 246:322
$ws.Rows(245).Insert()

# Populate the new row 245 with the new data record
$ws.Range("A245").Value2 = 10
$ws.Range("B245").Value2 = "Vega Modelo de Temuco"
$ws.Range("C245").Value2 = "La Araucanía"
$ws.Range("D245").Value2 = 44985
$ws.Range("E245").Value2 = 9
$ws.Range("F245").Value2 = "Fruta"
$ws.Range("G245").Value2 = 100103
$ws.Range("H245").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I245").Value2 = 100103002
$ws.Range("J245").Value2 = "Ciruela"
$ws.Range("K245").Value2 = "Blue Giant"
$ws.Range("L245").Value2 = "Primera"
$ws.Range("M245").Value2 = 145
$ws.Range("N245").Value2 = 14000
$ws.Range("O245").Value2 = 15000
$ws.Range("P245").Value2 = 14414
$ws.Range("Q245").Value2 = '$/bandeja 18 kilos granel'
$ws.Range("R245").Value2 = "Región de O'Higgins"
$ws.Range("S245").Value2 = 801
$ws.Range("T245").Value2 = 18
